# Insert a new row for "KIA Soul New" above the existing "KIA Sportage" row
# (row 136), shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(136).Insert()

$ws.Cells.Item(136, 1).NumberFormat = "@"
$ws.Cells.Item(136, 1).Value = "372"
$ws.Cells.Item(136, 2).Value = "KIA"
$ws.Cells.Item(136, 3).Value = "Soul New"
$ws.Cells.Item(136, 4).Value = 1887000
$ws.Cells.Item(136, 5).Value = "https://saratov-avtohous.ru/katalog/kia/soul/3-restyling"
$ws.Cells.Item(136, 6).Value = 1887000
$ws.Cells.Item(136, 7).Value = "https://saratov-avtohous.ru/katalog/kia/soul/3-restyling"
